$wb = $excel.ActiveWorkbook

# --- Logs sheet: append a new row (row 18) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(18, 1).Value = "Demo inplannen"
$logs.Cells.Item(18, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(18, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(18, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(18, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(18, 6).Value = "2025-08-13 21:47:40"
$logs.Cells.Item(18, 7).Value = "Nee"
$logs.Cells.Item(18, 8).Value = "Ja"
$logs.Cells.Item(18, 9).Value = "Nee"
$logs.Cells.Item(18, 10).Value = "Nee"

# --- Logs sheet: extend the conditional-formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`17")
    $newRange = $logs.Range("$col`2:$col`18")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the count in B2 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 17
